$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table data: Nomor Meja, Kapasitas, Status
# All rows now set to "tersedia" and two new rows (10, 11) appended with
# "Nomor Meja" renumbered 1..11 (previously 2..9).
$data = @(
    @(1, 5, "tersedia"),
    @(2, 3, "tersedia"),
    @(3, 5, "tersedia"),
    @(4, 2, "tersedia"),
    @(5, 2, "tersedia"),
    @(6, 5, "tersedia"),
    @(7, 5, "tersedia"),
    @(8, 5, "tersedia"),
    @(9, 5, "tersedia"),
    @(10, 3, "tersedia"),
    @(11, 5, "tersedia")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    if ($row -gt 9) {
        # New rows: carry over the same formatting (borders, style) used by
        # the existing data rows before filling in the values.
        $ws.Range("A9:C9").Copy()
        $ws.Range("A$row`:C$row").PasteSpecial(-4122)
    }
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("D7").Select()
